$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("D2893_2_bg")
$ws2.Range("C10").Value = 0.349
$ws2.Range("D10").Value = 0.348
$ws2.Range("E10").Value = 0.35
$ws2.Range("H10").Value = 0.349
$ws2.Range("I10").Value = 0.348
$ws2.Range("L10").Value = 0.348
$ws2.Range("M10").Value = 0.351
$ws2.Range("N10").Value = 0.349
$ws2.Range("O10").Value = 0.349
$ws2.Range("Q10").Value = 0.35
$ws2.Range("S10").Value = 0.354
$ws2.Range("W10").Value = 0.35
$ws2.Range("Y10").Value = 0.348
$ws2.Range("Z10").Value = 0.354
$ws2.Range("B11").Value = 0.06900000000000001
$ws2.Range("C11").Value = 0.068
$ws2.Range("D11").Value = 0.076
$ws2.Range("E11").Value = 0.062
$ws2.Range("F11").Value = 0.064
$ws2.Range("G11").Value = 0.037
$ws2.Range("H11").Value = 0.07199999999999999
$ws2.Range("I11").Value = 0.079
$ws2.Range("J11").Value = 0.053
$ws2.Range("K11").Value = 0.05
$ws2.Range("L11").Value = 0.075
$ws2.Range("M11").Value = 0.05
$ws2.Range("N11").Value = 0.07099999999999999
$ws2.Range("O11").Value = 0.07000000000000001
$ws2.Range("P11").Value = 0.066
$ws2.Range("Q11").Value = 0.059
$ws2.Range("R11").Value = 0.052
$ws2.Range("S11").Value = 0.025
$ws2.Range("T11").Value = 0.042
$ws2.Range("U11").Value = 0.027
$ws2.Range("V11").Value = 0.061
$ws2.Range("W11").Value = 0.059
$ws2.Range("X11").Value = 0.016
$ws2.Range("Y11").Value = 0.025
$ws2.Range("Z11").Value = 0.079
$ws2.Range("B16").Value = 45.631
$ws2.Range("C16").Value = 45.632
$ws2.Range("D16").Value = 45.624
$ws2.Range("E16").Value = 45.637
$ws2.Range("F16").Value = 45.635
$ws2.Range("G16").Value = 45.66
$ws2.Range("H16").Value = 45.628
$ws2.Range("I16").Value = 45.622
$ws2.Range("J16").Value = 45.645
$ws2.Range("K16").Value = 45.648
$ws2.Range("L16").Value = 45.625
$ws2.Range("M16").Value = 45.648
$ws2.Range("N16").Value = 45.629
$ws2.Range("O16").Value = 45.63
$ws2.Range("P16").Value = 45.634
$ws2.Range("Q16").Value = 45.64
$ws2.Range("R16").Value = 45.646
$ws2.Range("S16").Value = 45.67
$ws2.Range("T16").Value = 45.655
$ws2.Range("U16").Value = 45.668
$ws2.Range("V16").Value = 45.638
$ws2.Range("W16").Value = 45.64
$ws2.Range("X16").Value = 0.014
$ws2.Range("Y16").Value = 45.622
$ws2.Range("Z16").Value = 45.67
$ws2.Range("B28").Value = 3.12
$ws2.Range("C28").Value = 3.121
$ws2.Range("D28").Value = 3.113
$ws2.Range("E28").Value = 3.128
$ws2.Range("F28").Value = 3.126
$ws2.Range("G28").Value = 3.153
$ws2.Range("H28").Value = 3.118
$ws2.Range("I28").Value = 3.111
$ws2.Range("J28").Value = 3.136
$ws2.Range("K28").Value = 3.139
$ws2.Range("L28").Value = 3.114
$ws2.Range("M28").Value = 3.14
$ws2.Range("N28").Value = 3.118
$ws2.Range("O28").Value = 3.119
$ws2.Range("P28").Value = 3.123
$ws2.Range("Q28").Value = 3.13
$ws2.Range("R28").Value = 3.137
$ws2.Range("S28").Value = 3.164
$ws2.Range("T28").Value = 3.147
$ws2.Range("U28").Value = 3.163
$ws2.Range("V28").Value = 3.128
$ws2.Range("W28").Value = 3.131
$ws2.Range("X28").Value = 0.016
$ws2.Range("Y28").Value = 3.111
$ws2.Range("Z28").Value = 3.164
$ws2.Range("B29").Value = 0.06900000000000001
$ws2.Range("C29").Value = 0.068
$ws2.Range("D29").Value = 0.076
$ws2.Range("E29").Value = 0.062
$ws2.Range("F29").Value = 0.064
$ws2.Range("G29").Value = 0.037
$ws2.Range("H29").Value = 0.07199999999999999
$ws2.Range("I29").Value = 0.079
$ws2.Range("J29").Value = 0.053
$ws2.Range("K29").Value = 0.05
$ws2.Range("L29").Value = 0.075
$ws2.Range("M29").Value = 0.05
$ws2.Range("N29").Value = 0.07099999999999999
$ws2.Range("O29").Value = 0.07000000000000001
$ws2.Range("P29").Value = 0.066
$ws2.Range("Q29").Value = 0.059
$ws2.Range("R29").Value = 0.052
$ws2.Range("S29").Value = 0.025
$ws2.Range("T29").Value = 0.042
$ws2.Range("U29").Value = 0.027
$ws2.Range("V29").Value = 0.061
$ws2.Range("W29").Value = 0.059
$ws2.Range("X29").Value = 0.016
$ws2.Range("Y29").Value = 0.025
$ws2.Range("Z29").Value = 0.079
$ws2.Range("B37").Value = 15.927
$ws2.Range("C37").Value = 15.927
$ws2.Range("D37").Value = 15.929
$ws2.Range("E37").Value = 15.925
$ws2.Range("F37").Value = 15.925
$ws2.Range("H37").Value = 15.928
$ws2.Range("I37").Value = 15.93
$ws2.Range("J37").Value = 15.922
$ws2.Range("L37").Value = 15.929
$ws2.Range("M37").Value = 15.921
$ws2.Range("N37").Value = 15.928
$ws2.Range("O37").Value = 15.927
$ws2.Range("P37").Value = 15.926
$ws2.Range("Q37").Value = 15.924
$ws2.Range("R37").Value = 15.922
$ws2.Range("S37").Value = 15.914
$ws2.Range("T37").Value = 15.919
$ws2.Range("W37").Value = 15.924
$ws2.Range("X37").Value = 0.005
$ws2.Range("Y37").Value = 15.914
$ws2.Range("Z37").Value = 15.93
$ws2.Range("B38").Value = 5.893
$ws2.Range("E38").Value = 5.892
$ws2.Range("H38").Value = 5.893
$ws2.Range("I38").Value = 5.894
$ws2.Range("J38").Value = 5.891
$ws2.Range("N38").Value = 5.893
$ws2.Range("O38").Value = 5.893
$ws2.Range("Q38").Value = 5.892
$ws2.Range("R38").Value = 5.891
$ws2.Range("S38").Value = 5.888
$ws2.Range("Y38").Value = 5.888
$ws2.Range("Z38").Value = 5.894
$ws2.Range("C39").Value = 5.112
$ws2.Range("E39").Value = 5.111
$ws2.Range("H39").Value = 5.112
$ws2.Range("I39").Value = 5.113
$ws2.Range("M39").Value = 5.11
$ws2.Range("N39").Value = 5.112
$ws2.Range("O39").Value = 5.112
$ws2.Range("Q39").Value = 5.111
$ws2.Range("U39").Value = 5.108
$ws2.Range("W39").Value = 5.111
$ws2.Range("X39").Value = 0.002
$ws2.Range("Z39").Value = 5.113
$ws2.Range("D40").Value = 4.339
$ws2.Range("F40").Value = 4.338
$ws2.Range("I40").Value = 4.339
$ws2.Range("J40").Value = 4.337
$ws2.Range("L40").Value = 4.339
$ws2.Range("P40").Value = 4.338
$ws2.Range("R40").Value = 4.337
$ws2.Range("U40").Value = 4.335
$ws2.Range("V40").Value = 4.337
$ws2.Range("Z40").Value = 4.339
$ws2.Range("I41").Value = 2.083
$ws2.Range("J41").Value = 2.082
$ws2.Range("K41").Value = 2.082
$ws2.Range("M41").Value = 2.082
$ws2.Range("R41").Value = 2.082
$ws2.Range("S41").Value = 2.081
$ws2.Range("Y41").Value = 2.081
$ws2.Range("Z41").Value = 2.083
$ws2.Range("B45").Value = 6.998
$ws2.Range("C45").Value = 7
$ws2.Range("D45").Value = 6.983
$ws2.Range("E45").Value = 7.014
$ws2.Range("F45").Value = 7.009
$ws2.Range("G45").Value = 7.067
$ws2.Range("H45").Value = 6.992
$ws2.Range("I45").Value = 6.978
$ws2.Range("J45").Value = 7.032
$ws2.Range("K45").Value = 7.038
$ws2.Range("L45").Value = 6.985
$ws2.Range("M45").Value = 7.039
$ws2.Range("N45").Value = 6.993
$ws2.Range("O45").Value = 6.996
$ws2.Range("P45").Value = 7.004
$ws2.Range("Q45").Value = 7.019
$ws2.Range("R45").Value = 7.034
$ws2.Range("S45").Value = 7.09
$ws2.Range("T45").Value = 7.055
$ws2.Range("U45").Value = 7.087
$ws2.Range("V45").Value = 7.014
$ws2.Range("W45").Value = 7.02
$ws2.Range("X45").Value = 0.033
$ws2.Range("Y45").Value = 6.978
$ws2.Range("Z45").Value = 7.09
$ws2.Range("B46").Value = 0.1
$ws2.Range("C46").Value = 0.098
$ws2.Range("D46").Value = 0.11
$ws2.Range("E46").Value = 0.089
$ws2.Range("F46").Value = 0.092
$ws2.Range("G46").Value = 0.053
$ws2.Range("H46").Value = 0.104
$ws2.Range("I46").Value = 0.114
$ws2.Range("J46").Value = 0.077
$ws2.Range("K46").Value = 0.07199999999999999
$ws2.Range("L46").Value = 0.109
$ws2.Range("M46").Value = 0.07199999999999999
$ws2.Range("N46").Value = 0.103
$ws2.Range("O46").Value = 0.101
$ws2.Range("P46").Value = 0.095
$ws2.Range("Q46").Value = 0.08500000000000001
$ws2.Range("R46").Value = 0.075
$ws2.Range("S46").Value = 0.036
$ws2.Range("T46").Value = 0.061
$ws2.Range("U46").Value = 0.039
$ws2.Range("V46").Value = 0.08799999999999999
$ws2.Range("W46").Value = 0.08400000000000001
$ws2.Range("X46").Value = 0.022
$ws2.Range("Y46").Value = 0.036
$ws2.Range("Z46").Value = 0.114
$ws2.Range("J51").Value = 57.616
$ws2.Range("R51").Value = 57.616
$ws2.Range("W51").Value = 57.616

$ws3 = $wb.Worksheets.Item("D2893_3_bg_apf")
$ws3.Range("F10").Value = 0.348
$ws3.Range("H10").Value = 0.347
$ws3.Range("I10").Value = 0.346
$ws3.Range("L10").Value = 0.347
$ws3.Range("M10").Value = 0.35
$ws3.Range("N10").Value = 0.347
$ws3.Range("O10").Value = 0.347
$ws3.Range("Q10").Value = 0.349
$ws3.Range("S10").Value = 0.353
$ws3.Range("U10").Value = 0.353
$ws3.Range("V10").Value = 0.349
$ws3.Range("Y10").Value = 0.346
$ws3.Range("Z10").Value = 0.353
$ws3.Range("B11").Value = 0.083
$ws3.Range("C11").Value = 0.082
$ws3.Range("D11").Value = 0.092
$ws3.Range("E11").Value = 0.074
$ws3.Range("F11").Value = 0.077
$ws3.Range("G11").Value = 0.044
$ws3.Range("H11").Value = 0.08699999999999999
$ws3.Range("I11").Value = 0.095
$ws3.Range("J11").Value = 0.064
$ws3.Range("K11").Value = 0.061
$ws3.Range("L11").Value = 0.091
$ws3.Range("M11").Value = 0.06
$ws3.Range("N11").Value = 0.08599999999999999
$ws3.Range("O11").Value = 0.08500000000000001
$ws3.Range("P11").Value = 0.08
$ws3.Range("Q11").Value = 0.07099999999999999
$ws3.Range("R11").Value = 0.063
$ws3.Range("S11").Value = 0.03
$ws3.Range("T11").Value = 0.051
$ws3.Range("U11").Value = 0.032
$ws3.Range("V11").Value = 0.074
$ws3.Range("W11").Value = 0.07099999999999999
$ws3.Range("X11").Value = 0.019
$ws3.Range("Y11").Value = 0.03
$ws3.Range("Z11").Value = 0.095
$ws3.Range("B16").Value = 45.618
$ws3.Range("C16").Value = 45.619
$ws3.Range("D16").Value = 45.611
$ws3.Range("E16").Value = 45.626
$ws3.Range("F16").Value = 45.624
$ws3.Range("G16").Value = 45.653
$ws3.Range("H16").Value = 45.615
$ws3.Range("I16").Value = 45.608
$ws3.Range("J16").Value = 45.635
$ws3.Range("K16").Value = 45.638
$ws3.Range("L16").Value = 45.612
$ws3.Range("M16").Value = 45.639
$ws3.Range("N16").Value = 45.616
$ws3.Range("O16").Value = 45.617
$ws3.Range("P16").Value = 45.622
$ws3.Range("Q16").Value = 45.629
$ws3.Range("R16").Value = 45.636
$ws3.Range("S16").Value = 45.665
$ws3.Range("T16").Value = 45.647
$ws3.Range("U16").Value = 45.663
$ws3.Range("V16").Value = 45.626
$ws3.Range("W16").Value = 45.629
$ws3.Range("X16").Value = 0.017
$ws3.Range("Y16").Value = 45.608
$ws3.Range("Z16").Value = 45.665
$ws3.Range("B28").Value = 3.106
$ws3.Range("C28").Value = 3.107
$ws3.Range("D28").Value = 3.098
$ws3.Range("E28").Value = 3.115
$ws3.Range("F28").Value = 3.113
$ws3.Range("G28").Value = 3.145
$ws3.Range("H28").Value = 3.103
$ws3.Range("I28").Value = 3.095
$ws3.Range("J28").Value = 3.125
$ws3.Range("K28").Value = 3.129
$ws3.Range("L28").Value = 3.099
$ws3.Range("M28").Value = 3.129
$ws3.Range("N28").Value = 3.104
$ws3.Range("O28").Value = 3.105
$ws3.Range("P28").Value = 3.11
$ws3.Range("Q28").Value = 3.118
$ws3.Range("R28").Value = 3.127
$ws3.Range("S28").Value = 3.159
$ws3.Range("T28").Value = 3.139
$ws3.Range("U28").Value = 3.157
$ws3.Range("V28").Value = 3.116
$ws3.Range("W28").Value = 3.119
$ws3.Range("X28").Value = 0.019
$ws3.Range("Y28").Value = 3.095
$ws3.Range("Z28").Value = 3.159
$ws3.Range("B29").Value = 0.083
$ws3.Range("C29").Value = 0.082
$ws3.Range("D29").Value = 0.092
$ws3.Range("E29").Value = 0.074
$ws3.Range("F29").Value = 0.077
$ws3.Range("G29").Value = 0.044
$ws3.Range("H29").Value = 0.08699999999999999
$ws3.Range("I29").Value = 0.095
$ws3.Range("J29").Value = 0.064
$ws3.Range("K29").Value = 0.061
$ws3.Range("L29").Value = 0.091
$ws3.Range("M29").Value = 0.06
$ws3.Range("N29").Value = 0.08599999999999999
$ws3.Range("O29").Value = 0.08500000000000001
$ws3.Range("P29").Value = 0.08
$ws3.Range("Q29").Value = 0.07099999999999999
$ws3.Range("R29").Value = 0.063
$ws3.Range("S29").Value = 0.03
$ws3.Range("T29").Value = 0.051
$ws3.Range("U29").Value = 0.032
$ws3.Range("V29").Value = 0.074
$ws3.Range("W29").Value = 0.07099999999999999
$ws3.Range("X29").Value = 0.019
$ws3.Range("Y29").Value = 0.03
$ws3.Range("Z29").Value = 0.095
$ws3.Range("C37").Value = 15.931
$ws3.Range("D37").Value = 15.934
$ws3.Range("E37").Value = 15.929
$ws3.Range("F37").Value = 15.929
$ws3.Range("H37").Value = 15.932
$ws3.Range("I37").Value = 15.935
$ws3.Range("J37").Value = 15.926
$ws3.Range("L37").Value = 15.934
$ws3.Range("M37").Value = 15.924
$ws3.Range("N37").Value = 15.932
$ws3.Range("O37").Value = 15.932
$ws3.Range("P37").Value = 15.93
$ws3.Range("Q37").Value = 15.928
$ws3.Range("R37").Value = 15.925
$ws3.Range("S37").Value = 15.915
$ws3.Range("T37").Value = 15.921
$ws3.Range("U37").Value = 15.916
$ws3.Range("W37").Value = 15.927
$ws3.Range("X37").Value = 0.006
$ws3.Range("Y37").Value = 15.915
$ws3.Range("Z37").Value = 15.935
$ws3.Range("E38").Value = 5.893
$ws3.Range("F38").Value = 5.894
$ws3.Range("H38").Value = 5.895
$ws3.Range("I38").Value = 5.896
$ws3.Range("J38").Value = 5.892
$ws3.Range("K38").Value = 5.892
$ws3.Range("L38").Value = 5.895
$ws3.Range("M38").Value = 5.892
$ws3.Range("O38").Value = 5.894
$ws3.Range("P38").Value = 5.894
$ws3.Range("Q38").Value = 5.893
$ws3.Range("R38").Value = 5.892
$ws3.Range("W38").Value = 5.893
$ws3.Range("Z38").Value = 5.896
$ws3.Range("C39").Value = 5.113
$ws3.Range("D39").Value = 5.114
$ws3.Range("I39").Value = 5.114
$ws3.Range("L39").Value = 5.114
$ws3.Range("M39").Value = 5.111
$ws3.Range("O39").Value = 5.113
$ws3.Range("P39").Value = 5.113
$ws3.Range("Q39").Value = 5.112
$ws3.Range("S39").Value = 5.108
$ws3.Range("W39").Value = 5.112
$ws3.Range("Y39").Value = 5.108
$ws3.Range("Z39").Value = 5.114
$ws3.Range("F40").Value = 4.339
$ws3.Range("H40").Value = 4.34
$ws3.Range("J40").Value = 4.338
$ws3.Range("L40").Value = 4.34
$ws3.Range("N40").Value = 4.34
$ws3.Range("R40").Value = 4.338
$ws3.Range("S40").Value = 4.335
$ws3.Range("V40").Value = 4.338
$ws3.Range("X40").Value = 0.002
$ws3.Range("Y40").Value = 4.335
$ws3.Range("O41").Value = 2.083
$ws3.Range("P41").Value = 2.083
$ws3.Range("D44").Value = 1.109
$ws3.Range("I44").Value = 1.109
$ws3.Range("L44").Value = 1.109
$ws3.Range("Z44").Value = 1.109
$ws3.Range("B45").Value = 6.968
$ws3.Range("C45").Value = 6.97
$ws3.Range("D45").Value = 6.95
$ws3.Range("E45").Value = 6.987
$ws3.Range("F45").Value = 6.982
$ws3.Range("G45").Value = 7.051
$ws3.Range("H45").Value = 6.961
$ws3.Range("I45").Value = 6.944
$ws3.Range("J45").Value = 7.008
$ws3.Range("K45").Value = 7.016
$ws3.Range("L45").Value = 6.953
$ws3.Range("M45").Value = 7.017
$ws3.Range("N45").Value = 6.962
$ws3.Range("O45").Value = 6.966
$ws3.Range("P45").Value = 6.976
$ws3.Range("Q45").Value = 6.993
$ws3.Range("R45").Value = 7.012
$ws3.Range("S45").Value = 7.079
$ws3.Range("T45").Value = 7.037
$ws3.Range("U45").Value = 7.075
$ws3.Range("V45").Value = 6.988
$ws3.Range("W45").Value = 6.995
$ws3.Range("X45").Value = 0.04
$ws3.Range("Y45").Value = 6.944
$ws3.Range("Z45").Value = 7.079
$ws3.Range("B46").Value = 0.12
$ws3.Range("C46").Value = 0.119
$ws3.Range("D46").Value = 0.132
$ws3.Range("E46").Value = 0.107
$ws3.Range("F46").Value = 0.111
$ws3.Range("G46").Value = 0.064
$ws3.Range("H46").Value = 0.125
$ws3.Range("I46").Value = 0.137
$ws3.Range("J46").Value = 0.092
$ws3.Range("K46").Value = 0.08699999999999999
$ws3.Range("L46").Value = 0.131
$ws3.Range("M46").Value = 0.08699999999999999
$ws3.Range("N46").Value = 0.124
$ws3.Range("O46").Value = 0.122
$ws3.Range("P46").Value = 0.115
$ws3.Range("Q46").Value = 0.103
$ws3.Range("R46").Value = 0.09
$ws3.Range("S46").Value = 0.044
$ws3.Range("T46").Value = 0.073
$ws3.Range("U46").Value = 0.047
$ws3.Range("V46").Value = 0.107
$ws3.Range("W46").Value = 0.102
$ws3.Range("X46").Value = 0.027
$ws3.Range("Y46").Value = 0.044
$ws3.Range("Z46").Value = 0.137
$ws3.Range("H48").Value = 0.579
$ws3.Range("L48").Value = 0.579
$ws3.Range("N48").Value = 0.579
$ws3.Range("O48").Value = 0.579
$ws3.Range("T51").Value = 57.617
